$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Rework the "RNA techniques" sentence in the Skills to be enhanced
#    paragraph:
#    "I will use my connections with the Krainer lab to gain training
#    in RNA techniques." ->
#    "I will also gain training in RNA techniques from the members of
#    the Krainer lab. Finally, during this project, I will train my
#    skill at grant writing by applying for a K99 grant during the
#    third year of my posdoc. "
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "I will use my connections with the ", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "I will also gain training in RNA techniques from the members of the ",
    2)

$null = $d.Content.Find.Execute(
    " lab to gain training in RNA techniques.", $true, $false, $false,
    $false, $false, $true, 1, $false,
    " lab. Finally, during this project, I will train my skill at grant writing by applying for a K99 grant during the third year of my posdoc. ",
    2)

# ------------------------------------------------------------------
# 2. Extend the paragraph that ends in "...prepare for my job search."
#    with the new teaching-TA sentence.
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "ob search.", $true, $false, $false, $false, $false, $true, 1, $false,
    "ob search. Finally, In order to strengthen my skill teaching I will TA the graduate level quantitative biology course at CSHL.",
    2)

# ------------------------------------------------------------------
# 3. Insert the new closing paragraph right after that paragraph.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("I will TA the graduate level quantitative biology course at CSHL.")) {
        $newPara = $p.Range.InsertParagraphAfter()
        $p2 = $d.Paragraphs.Item($i + 1)
        $p2.Range.InsertAfter("Together, these plans will prepare me for a career as independent researcher by helping me to learn applications of my expertise to new fields, deepening my knowledge of computational methods, building my network of mentors and potential collaborators, and strengthening my skills at essential nonscience tasks like teaching and grant writing.")
        break
    }
}

# ------------------------------------------------------------------
# 4. Remove the three single, underlined sub-heading paragraphs
#    ("Overall training goals", "Skills to be enhanced",
#    "Preparation for career plans") -- the sections now read as
#    continuous prose.
# ------------------------------------------------------------------
$headings = @("Overall training goals", "Skills to be enhanced", "Preparation for career plans")
foreach ($heading in $headings) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($heading)) {
            $p.Range.Delete()
            break
        }
    }
}
